# Update details on CV (experience "where" column + skills tech stack names)

$wb = $excel.ActiveWorkbook

# --- "experience" sheet: add/update the "where" (column C) values ---
$wsExp = $wb.Worksheets.Item("experience")

# Self-employed / Marine ecology consultant rows: add "Remote" as the work location
$wsExp.Range("C2").Value = "Remote"
$wsExp.Range("C3").Value = "Remote"
$wsExp.Range("C4").Value = "Remote"
$wsExp.Range("C5").Value = "Remote"
$wsExp.Range("C6").Value = "Remote"

# Project ecologist & Data analyst @ HaMaarag rows: location becomes hybrid
$wsExp.Range("C7").Value = "Tel Aviv, Israel/Hybrid"
$wsExp.Range("C8").Value = "Tel Aviv, Israel/Hybrid"
$wsExp.Range("C9").Value = "Tel Aviv, Israel/Hybrid"
$wsExp.Range("C10").Value = "Tel Aviv, Israel/Hybrid"

# Database administrator @ Sharks in Israel rows: location becomes remote
$wsExp.Range("C14").Value = "Israel/Remote"
$wsExp.Range("C15").Value = "Israel/Remote"
$wsExp.Range("C16").Value = "Israel/Remote"

# MSc research student @ Belmaker Lab rows: location becomes hybrid
$wsExp.Range("C21").Value = "Tel Aviv, Israel/Hybrid"
$wsExp.Range("C22").Value = "Tel Aviv, Israel/Hybrid"

# --- "skills" sheet: rename a couple of tech-stack entries ---
$wsSkills = $wb.Worksheets.Item("skills")

$wsSkills.Range("C15").Value = "Google Workspace"
$wsSkills.Range("C2").Value = "R Programming Language"

[void]$wsSkills.Range("C3").Select()

# Leave the "experience" sheet as the active tab/selection, matching the
# workbook's final UI state (cursor resting on C18).
[void]$wsExp.Activate()
[void]$wsExp.Range("C18").Select()
